# SAP Test Planner Healthcheck - "UI Validations added V.23"
# Adds a new test-case row (row 28) to the TestCases sheet:
#   SuiteName=SAP Regression Automation, Runmode=Yes, TC_id=NC_OP_24,
#   TestRail_id=5400136, Test Name=Verify UI Change Data Correction to Price Determined

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the last existing row (27) down into the new row (28) so the
#    new row inherits the same cell formatting/borders as the rest of the table.
$ws.Range("A27:E27").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)

# 2) Fill in the new test case's data (Test Name/TC_id added in this order so
#    the shared-string table grows in the same order as the source edit).
$ws.Range("A28").Value = "SAP Regression Automation"
$ws.Range("B28").Value = "Yes"
$ws.Range("E28").Value = "Verify UI Change Data Correction to Price Determined"
$ws.Range("C28").Value = "NC_OP_24"
$ws.Range("D28").Value = 5400136

# 3) Re-apply the common formatting used by the rest of the TC_id/TestRail_id
#    columns onto the existing rows 26:27 so their styling is consistent with
#    the rest of the table (matches how Excel normalized them on save).
$ws.Range("D25:E25").Copy()
$ws.Range("D26:E27").PasteSpecial(-4122)

# 4) Extend the Yes/No dropdown validation on column B to cover the new row.
[void]$ws.Range("B2:B28").Validation.Delete()
[void]$ws.Range("B2:B28").Validation.Add(3, 1, 1, '"Yes,No"')

# 5) Match the saved selection/cursor position.
[void]$ws.Range("E23").Select()
